$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Base Models")
$ws.Rows.AutoFit()
$ws.Columns.AutoFit()
Write-Output "done"
